$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Connections sheet: add "tax" / "no_tax" columns (D, E) and split each
# *_TAX row into a *_TAX_ON row (existing, renamed) and a new *_TAX_OFF
# row, so the tax can be toggled per scenario.
# ---------------------------------------------------------------------
$wsConn = $wb.Worksheets.Item("Connections")

# New yellow-filled helper columns D:E, rows 1-18
$wsConn.Range("D1:E18").Interior.Color = 65535

# Rename existing *_TAX rows to *_TAX_ON, move their "Y" flag into the new
# "tax" (D) column, and clear the old C column.
$wsConn.Range("A9").Value = "COAL_TAX_ON"
$wsConn.Range("C9").ClearContents()
$wsConn.Range("D9").Value = "Y"

$wsConn.Range("A10").Value = "DSL_TAX_ON"
$wsConn.Range("C10").ClearContents()
$wsConn.Range("D10").Value = "Y"

$wsConn.Range("A11").Value = "MSW_LF_TAX_ON"
$wsConn.Range("C11").ClearContents()
$wsConn.Range("D11").Value = "Y"

$wsConn.Range("A12").Value = "OIL_TAX_ON"
$wsConn.Range("C12").ClearContents()
$wsConn.Range("D12").Value = "Y"

$wsConn.Range("A13").Value = "NATGAS_TAX_ON"
$wsConn.Range("C13").ClearContents()
$wsConn.Range("D13").Value = "Y"

# New *_TAX_OFF rows: "Y" goes in C (B column, same as before) and in the
# new "no_tax" (E) column.
$wsConn.Range("A14").Value = "COAL_TAX_OFF"
$wsConn.Range("C14").Value = "Y"
$wsConn.Range("E14").Value = "Y"

$wsConn.Range("A15").Value = "DSL_TAX_OFF"
$wsConn.Range("C15").Value = "Y"
$wsConn.Range("E15").Value = "Y"

$wsConn.Range("A16").Value = "MSW_LF_TAX_OFF"
$wsConn.Range("C16").Value = "Y"
$wsConn.Range("E16").Value = "Y"

$wsConn.Range("A17").Value = "OIL_TAX_OFF"
$wsConn.Range("C17").Value = "Y"
$wsConn.Range("E17").Value = "Y"

$wsConn.Range("A18").Value = "NATGAS_TAX_OFF"
$wsConn.Range("C18").Value = "Y"
$wsConn.Range("E18").Value = "Y"

# Header row for the two new helper columns
$wsConn.Range("D1").Value = "tax"
$wsConn.Range("E1").Value = "no_tax"

# Widen column A now that it holds longer names like "NATGAS_TAX_OFF"
$wsConn.Columns.Item(1).ColumnWidth = 19.14

# ---------------------------------------------------------------------
# SolverSettings sheet used to be the active tab; move its selection
# without leaving it the active tab (Connections takes over below).
# ---------------------------------------------------------------------
$wsSolver = $wb.Worksheets.Item("SolverSettings")
$wsSolver.Activate()
$wsSolver.Range("H6").Select()

# Make Connections the active sheet/tab with a new selection -- this must
# be the last sheet activated so it ends up as the workbook's active tab.
$wsConn.Activate()
$wsConn.Range("H10").Select()
